$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f963ff84ed7d7a57a79365740cd4f8ea0b65349/e2e/fc8472b5-ff02-4bf0-ba73-d2bb1cf1efac.md"
$displayName = "fc8472b5-ff02-4bf0-ba73-d2bb1cf1efac.md"
$versionMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c4af0cb138841a6da3f7a688257e5b53887c8287/e2e/fc8472b5-ff02-4bf0-ba73-d2bb1cf1efac.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f963ff84ed7d7a57a79365740cd4f8ea0b65349/e2e/fc8472b5-ff02-4bf0-ba73-d2bb1cf1efac.md."

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestUrl, "", "", $displayName)
$wsZh.Range("I7").Style = "HyperLink"

$wsZh.Range("J7").Value = "fc8472b5-ff02-4bf0-ba73-d2bb1cf1efac.4692d72a51657069896a8d384d391c2ce0ad0148.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-18 04:55:47"
$wsZh.Range("P7").Value = $versionMessage

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestUrl, "", "", $displayName)
$wsDe.Range("I7").Style = "HyperLink"

$wsDe.Range("J7").Value = "fc8472b5-ff02-4bf0-ba73-d2bb1cf1efac.4692d72a51657069896a8d384d391c2ce0ad0148.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-18 04:55:54"
$wsDe.Range("P7").Value = $versionMessage
